$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.309.70"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "'1.876.41"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'0.7121"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").Value = "'242.34"
$ws.Range("E6").Value = "  +0.73%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.3104"
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("D9").Value = "'0.07765"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "'24.92"
$ws.Range("D11").Value = "'0.08504"
$ws.Range("E11").Value = "  +3.00%  "
$ws.Range("D12").Value = "'1.880.02"
$ws.Range("E12").Value = "  +1.70%  "
$ws.Range("D13").Value = "'5.212"
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").Value = "'91.41"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").Value = "'29.310.22"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").Value = "'0.000008232"
$ws.Range("E17").Value = "  +5.57%  "
$ws.Range("E18").Value = "  +2.41%  "
$ws.Range("D19").Value = "'242.65"
$ws.Range("E19").Value = "  -0.70%  "
$ws.Range("D20").Value = "'2.133.61"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("D22").Value = "'0.9998"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "'7.821"
$ws.Range("E23").Value = "  -2.03%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'0.1625"
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("D26").Value = "'162.76"
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").Value = "'9.026"
$ws.Range("E29").Value = "  +1.02%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "'4.326"
$ws.Range("E31").Value = "  +3.24%  "
$ws.Range("D32").Value = "'1.279"
$ws.Range("E32").Value = "  -2.73%  "
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("D34").Value = "'1.932"
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("D36").Value = "'0.7398"
$ws.Range("E36").Value = "  +1.66%  "
$ws.Range("D37").Value = "'2.686"
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("D38").Value = "'0.01865"
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("D39").Value = "'2.723"
$ws.Range("E39").Value = "  +1.40%  "
$ws.Range("D40").Value = "'1.174.36"
$ws.Range("E40").Value = "  +2.15%  "
$ws.Range("D41").Value = "'6.382"
$ws.Range("E41").Value = "  +4.09%  "
$ws.Range("D42").Value = "'0.8889"
$ws.Range("E42").Value = "  -1.88%  "
$ws.Range("D43").Value = "'72.96"
$ws.Range("E43").Value = "  +0.95%  "
$ws.Range("D44").Value = "'106.29"
$ws.Range("E44").Value = "  +4.45%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "'2.029.78"
$ws.Range("E46").Value = "  +1.05%  "
$ws.Range("E47").Value = "  +2.48%  "
$ws.Range("D48").Value = "'0.5207"
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("D49").Value = "'0.00000000121"
$ws.Range("E49").Value = "  +1.60%  "
$ws.Range("D50").Value = "'9.396"
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("D51").Value = "'0.4309"
$ws.Range("E51").Value = "  +1.08%  "
